$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "19+33=52"
$t.Cell(1, 2).Range.Text = "33-29=4"
$t.Cell(1, 3).Range.Text = "39+57=96"
$t.Cell(1, 4).Range.Text = "65-39=26"
$t.Cell(1, 5).Range.Text = "6+8=14"
$t.Cell(2, 1).Range.Text = "34-19=15"
$t.Cell(2, 2).Range.Text = "19+29=48"
$t.Cell(2, 3).Range.Text = "66-38=28"
$t.Cell(2, 4).Range.Text = "29+55=84"
$t.Cell(2, 5).Range.Text = "35+29=64"
$t.Cell(3, 1).Range.Text = "52-13=39"
$t.Cell(3, 2).Range.Text = "94-39=55"
$t.Cell(3, 3).Range.Text = "93-84=9"
$t.Cell(3, 4).Range.Text = "59+28=87"
$t.Cell(3, 5).Range.Text = "39+25=64"
$t.Cell(4, 1).Range.Text = "61-59=2"
$t.Cell(4, 2).Range.Text = "66-59=7"
$t.Cell(4, 3).Range.Text = "41-28=13"
$t.Cell(4, 4).Range.Text = "42-28=14"
$t.Cell(4, 5).Range.Text = "61-17=44"
$t.Cell(5, 1).Range.Text = "60-4=56"
$t.Cell(5, 2).Range.Text = "50-29=21"
$t.Cell(5, 3).Range.Text = "8+88=96"
$t.Cell(5, 4).Range.Text = "54-46=8"
$t.Cell(5, 5).Range.Text = "14+58=72"
$t.Cell(6, 1).Range.Text = "17+28=45"
$t.Cell(6, 2).Range.Text = "50-14=36"
$t.Cell(6, 3).Range.Text = "56+9=65"
$t.Cell(6, 4).Range.Text = "8+19=27"
$t.Cell(6, 5).Range.Text = "58+26=84"
$t.Cell(7, 1).Range.Text = "42-26=16"
$t.Cell(7, 2).Range.Text = "30-4=26"
$t.Cell(7, 3).Range.Text = "44+7=51"
$t.Cell(7, 4).Range.Text = "70-34=36"
$t.Cell(7, 5).Range.Text = "36+16=52"
$t.Cell(8, 1).Range.Text = "26+29=55"
$t.Cell(8, 2).Range.Text = "84-26=58"
$t.Cell(8, 3).Range.Text = "93-6=87"
$t.Cell(8, 4).Range.Text = "93-5=88"
$t.Cell(8, 5).Range.Text = "18+17=35"
$t.Cell(9, 1).Range.Text = "9+3=12"
$t.Cell(9, 2).Range.Text = "85-66=19"
$t.Cell(9, 3).Range.Text = "38+56=94"
$t.Cell(9, 4).Range.Text = "35+18=53"
$t.Cell(9, 5).Range.Text = "42-18=24"
$t.Cell(10, 1).Range.Text = "79+19=98"
$t.Cell(10, 2).Range.Text = "45-29=16"
$t.Cell(10, 3).Range.Text = "4+69=73"
$t.Cell(10, 4).Range.Text = "76-38=38"
$t.Cell(10, 5).Range.Text = "31-28=3"
$t.Cell(11, 1).Range.Text = "85-69=16"
$t.Cell(11, 2).Range.Text = "39+13=52"
$t.Cell(11, 3).Range.Text = "77-29=48"
$t.Cell(11, 4).Range.Text = "74-5=69"
$t.Cell(11, 5).Range.Text = "65+29=94"
$t.Cell(12, 1).Range.Text = "9+69=78"
$t.Cell(12, 2).Range.Text = "16+8=24"
$t.Cell(12, 3).Range.Text = "30-18=12"
$t.Cell(12, 4).Range.Text = "31-3=28"
$t.Cell(12, 5).Range.Text = "18+58=76"
$t.Cell(13, 1).Range.Text = "34+47=81"
$t.Cell(13, 2).Range.Text = "43+49=92"
$t.Cell(13, 3).Range.Text = "25+67=92"
$t.Cell(13, 4).Range.Text = "71-64=7"
$t.Cell(13, 5).Range.Text = "6+69=75"
$t.Cell(14, 1).Range.Text = "18+36=54"
$t.Cell(14, 2).Range.Text = "59+16=75"
$t.Cell(14, 3).Range.Text = "59+9=68"
$t.Cell(14, 4).Range.Text = "60-38=22"
$t.Cell(14, 5).Range.Text = "41-29=12"
$t.Cell(15, 1).Range.Text = "26+36=62"
$t.Cell(15, 2).Range.Text = "54+39=93"
$t.Cell(15, 3).Range.Text = "8+39=47"
$t.Cell(15, 4).Range.Text = "36-7=29"
$t.Cell(15, 5).Range.Text = "55+38=93"
$t.Cell(16, 1).Range.Text = "13+78=91"
$t.Cell(16, 2).Range.Text = "19+17=36"
$t.Cell(16, 3).Range.Text = "26+67=93"
$t.Cell(16, 4).Range.Text = "25+39=64"
$t.Cell(16, 5).Range.Text = "59+27=86"
$t.Cell(17, 1).Range.Text = "72-29=43"
$t.Cell(17, 2).Range.Text = "21-8=13"
$t.Cell(17, 3).Range.Text = "77+16=93"
$t.Cell(17, 4).Range.Text = "7+85=92"
$t.Cell(17, 5).Range.Text = "91-15=76"
$t.Cell(18, 1).Range.Text = "52-25=27"
$t.Cell(18, 2).Range.Text = "35+38=73"
$t.Cell(18, 3).Range.Text = "86-59=27"
$t.Cell(18, 4).Range.Text = "7+29=36"
$t.Cell(18, 5).Range.Text = "50-28=22"
$t.Cell(19, 1).Range.Text = "33-28=5"
$t.Cell(19, 2).Range.Text = "24+7=31"
$t.Cell(19, 3).Range.Text = "48-29=19"
$t.Cell(19, 4).Range.Text = "93-47=46"
$t.Cell(19, 5).Range.Text = "19+28=47"
$t.Cell(20, 1).Range.Text = "27+24=51"
$t.Cell(20, 2).Range.Text = "23-15=8"
$t.Cell(20, 3).Range.Text = "59+24=83"
$t.Cell(20, 4).Range.Text = "27+19=46"
$t.Cell(20, 5).Range.Text = "52-44=8"
